$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that are no longer present in the new data
# ("SC 92" at row 28 and "RM 232" at row 26). Delete the lower
# row first so the higher row index is unaffected.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# Apply the updated / newly-missing values (row numbers below are
# the FINAL row numbers, after the two rows above were removed).

# RM 14 (row 5): E now missing
$ws.Range("E5").ClearContents()

# RM 58 (row 11): E now has a value
$ws.Range("E11").Value = -7.9

# RM 125 (row 19): D now has a value, E now missing
$ws.Range("D19").Value = -15.5
$ws.Range("E19").ClearContents()

# RM 135 (row 21): D now missing
$ws.Range("D21").ClearContents()

# RM 140 (row 23): D and E now have values
$ws.Range("D23").Value = -13.9
$ws.Range("E23").Value = -7

# RM 145 (row 25): E now has a value
$ws.Range("E25").Value = -7.1

# SC 5 (row 26): C now missing
$ws.Range("C26").ClearContents()

# SC 101 (row 27): C now has a value, D and E now missing
$ws.Range("C27").Value = 10
$ws.Range("D27").ClearContents()
$ws.Range("E27").ClearContents()

# SC 119 (row 29): C and E now missing
$ws.Range("C29").ClearContents()
$ws.Range("E29").ClearContents()

# SC 232 (row 33): D and E now have values
$ws.Range("D33").Value = -14.1
$ws.Range("E33").Value = -10.7
